# DTT-Assessment-Hour-Log.xlsx — log the last week's hours.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 — tweak the description text (added a trailing period / reworded).
$ws.Range("D5").Value = "Started creating all the files needed for the project and started working with components."

# Row 6
$ws.Range("A6").Value = "Continue working on houseList and filterComponent"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 45242
$ws.Range("D6").Value = "Continued working on filterComponent, houseList and homeView to make evertything work as expected."

# Row 7
$ws.Range("A7").Value = "Continue working on houseList and HomeView"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 45245
$ws.Range("D7").Value = "Fixing errors that were causing things to render incorrectly. Also added button for sorting."

# Row 8
$ws.Range("A8").Value = "Working on sorting and filtering"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 45246
$ws.Range("D8").Value = "Fixing some bugs that caused the site to be unable to sort or filter houses based on the users input."

# Row 9
$ws.Range("A9").Value = "Adding create new button and page"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 45247
$ws.Range("D9").Value = "Adding button for creating new house, also added the page where you can create that new house. "

# Widen column A and D so the longer log text is readable.
$ws.Columns("A").ColumnWidth = 27.75
$ws.Columns("D").ColumnWidth = 51.625

# Scroll / zoom / selection to match where the author left the view.
$win = $excel.ActiveWindow
$win.Zoom = 172
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("B10").Select()

# B30's SUMIF total recalculates automatically from the new hour entries.
